$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 54) with the next forecast data point.
# Copy the style of the prior row's date cell (A53) so A54 gets the
# same date number-format/style as the rest of column A, without
# introducing any new style definitions.
$row = 54

$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = [DateTime]::FromOADate(45986)
$ws.Cells.Item($row, 2).Value = 2025
$ws.Cells.Item($row, 3).Value = -0.08656168856399082
$ws.Cells.Item($row, 4).Value = 2026
$ws.Cells.Item($row, 5).Value = -0.02867614772544824

$wb.Save()
